# Update countries & provincias Spain
#
# The "Pais" sheet is a COVID-19 country leaderboard sorted (descending)
# by total cases (column B). A handful of countries received new case
# counts in this data refresh; because the sheet stays sorted by total
# cases, several rows whose counts are tied (or which changed rank)
# need their country name (column A) reassigned to match the new
# sort order, in addition to the numeric columns (B:H) being updated
# for the rows whose underlying stats actually changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp footer
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Marzo de 2020 a las 05:46"

$ws.Cells.Item(6, 2).Value = 26867
$ws.Cells.Item(6, 3).Value = 2660
$ws.Cells.Item(6, 5).Value = 26341

$ws.Cells.Item(47, 4).Value = 24
$ws.Cells.Item(47, 5).Value = 303

$ws.Cells.Item(60, 1).Value = "Colombia"
$ws.Cells.Item(60, 2).Value = 210
$ws.Cells.Item(60, 3).Value = 14
$ws.Cells.Item(60, 4).Value = 3
$ws.Cells.Item(60, 5).Value = 206
$ws.Cells.Item(60, 7).Value = 1

$ws.Cells.Item(61, 1).Value = "Croacia"
$ws.Cells.Item(61, 2).Value = 206
$ws.Cells.Item(61, 4).Value = 5
$ws.Cells.Item(61, 5).Value = 200
$ws.Cells.Item(61, 8).Value = 1

$ws.Cells.Item(110, 2).Value = 26
$ws.Cells.Item(110, 3).Value = 2
$ws.Cells.Item(110, 5).Value = 26

$ws.Cells.Item(116, 1).Value = "Ghana"

$ws.Cells.Item(118, 1).Value = "Puerto Rico"

$ws.Cells.Item(131, 1).Value = "Kirguistan"

$ws.Cells.Item(132, 1).Value = "Barbados"

$ws.Cells.Item(133, 1).Value = "Mauricio"
$ws.Cells.Item(133, 4).Value = 0
$ws.Cells.Item(133, 8).Value = 1

$ws.Cells.Item(134, 1).Value = "Costa de Marfil"
$ws.Cells.Item(134, 4).Value = 1
$ws.Cells.Item(134, 8).Value = 0

$ws.Cells.Item(140, 1).Value = "Aruba"
$ws.Cells.Item(140, 2).Value = 8
$ws.Cells.Item(140, 3).Value = 3
$ws.Cells.Item(140, 4).Value = 1

$ws.Cells.Item(142, 1).Value = "Kenia"
$ws.Cells.Item(142, 2).Value = 7
$ws.Cells.Item(142, 5).Value = 7

$ws.Cells.Item(143, 1).Value = "Guinea Ecuatorial"

$ws.Cells.Item(144, 1).Value = "Islas Virgenes de los Estados Unidos"

$ws.Cells.Item(145, 1).Value = "Tanzania"
$ws.Cells.Item(145, 2).Value = 6
$ws.Cells.Item(145, 5).Value = 6

$ws.Cells.Item(147, 1).Value = "Surinam"
$ws.Cells.Item(147, 5).Value = 5
$ws.Cells.Item(147, 8).Value = 0

$ws.Cells.Item(148, 1).Value = "Gabon"
$ws.Cells.Item(148, 4).Value = 0
$ws.Cells.Item(148, 8).Value = 1

$ws.Cells.Item(149, 1).Value = "Suazilandia"
$ws.Cells.Item(149, 3).Value = 3

$ws.Cells.Item(150, 1).Value = "Bahamas"

$ws.Cells.Item(151, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(151, 2).Value = 4
$ws.Cells.Item(151, 5).Value = 4

$ws.Cells.Item(152, 1).Value = "Zimbabue"

$ws.Cells.Item(153, 1).Value = "Namibia"

$ws.Cells.Item(154, 1).Value = "El Salvador"

$ws.Cells.Item(155, 1).Value = "San Bartolome"

$ws.Cells.Item(156, 1).Value = "Cabo Verde"

$ws.Cells.Item(157, 1).Value = "Congo"

$ws.Cells.Item(158, 1).Value = "Liberia"

$ws.Cells.Item(159, 1).Value = "Madagascar"

$ws.Cells.Item(160, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(160, 5).Value = 3
$ws.Cells.Item(160, 8).Value = 0

$ws.Cells.Item(161, 1).Value = "Islas Caimanes"

$ws.Cells.Item(162, 1).Value = "Curazao"
$ws.Cells.Item(162, 2).Value = 3
$ws.Cells.Item(162, 8).Value = 1

$ws.Cells.Item(163, 1).Value = "Haiti"

$ws.Cells.Item(164, 1).Value = "Butan"

$ws.Cells.Item(165, 1).Value = "Zambia"

$ws.Cells.Item(166, 1).Value = "Groenlandia"

$ws.Cells.Item(167, 1).Value = "Fiyi"

$ws.Cells.Item(168, 1).Value = "Santa Lucia"

$ws.Cells.Item(169, 1).Value = "Guinea"

$ws.Cells.Item(170, 1).Value = "Benin"

$ws.Cells.Item(171, 1).Value = "Bermudas"

$ws.Cells.Item(172, 1).Value = "Nicaragua"

$ws.Cells.Item(173, 1).Value = "Mauritania"

$ws.Cells.Item(174, 1).Value = "Isla de Man"

$ws.Cells.Item(175, 1).Value = "Angola"
$ws.Cells.Item(175, 5).Value = 2
$ws.Cells.Item(175, 8).Value = 0

$ws.Cells.Item(176, 1).Value = "Sudan"
$ws.Cells.Item(176, 2).Value = 2
$ws.Cells.Item(176, 8).Value = 1

$ws.Cells.Item(177, 1).Value = "San Martin (Parte Holandesa)"

$ws.Cells.Item(189, 1).Value = "Niger"

$ws.Cells.Item(190, 1).Value = "Timor Oriental"
